# Auto-generated edit script applying odds updates from the 2024-10-17 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.3
$ws.Range("J3").Value = 2.63
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 5.5
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 5
$ws.Range("Y3").Value = 9.5
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AH3").Value = 9
$ws.Range("AJ3").Value = 17
$ws.Range("AT3").Value = 2.25
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AX3").Value = 29
$ws.Range("BB3").Value = 451

# Row 4
$ws.Range("G4").Value = 1.73
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 1.95
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.38
$ws.Range("V4").Value = 1.53
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 6
$ws.Range("AH4").Value = 11
$ws.Range("AJ4").Value = 21
$ws.Range("AM4").Value = 67
$ws.Range("AO4").Value = 9.5
$ws.Range("AP4").Value = 26
$ws.Range("AR4").Value = 67
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.25
$ws.Range("AU4").Value = 10
$ws.Range("AW4").Value = 6.5
$ws.Range("BA4").Value = 201

# Row 5
$ws.Range("I5").Value = 2.4
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("W5").Value = 9
$ws.Range("AC5").Value = 10
$ws.Range("AI5").Value = 11
$ws.Range("AL5").Value = 19
$ws.Range("AN5").Value = 5
$ws.Range("AT5").Value = 2.75

# Row 6
$ws.Range("P6").Value = 4.05
$ws.Range("U6").Value = 1.9
$ws.Range("V6").Value = 1.86

# Row 9
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 2.25

# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 2.5
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 4.5
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 11
$ws.Range("O10").Value = 1.29
$ws.Range("P10").Value = 3.5
$ws.Range("Q10").Value = 1.98
$ws.Range("R10").Value = 1.88
$ws.Range("S10").Value = 1.4
$ws.Range("T10").Value = 2.75
$ws.Range("U10").Value = 1.83
$ws.Range("V10").Value = 1.83
$ws.Range("W10").Value = 7
$ws.Range("X10").Value = 8.5
$ws.Range("Y10").Value = 8.5
$ws.Range("Z10").Value = 15
$ws.Range("AA10").Value = 15
$ws.Range("AB10").Value = 26
$ws.Range("AC10").Value = 10
$ws.Range("AD10").Value = 7
$ws.Range("AE10").Value = 15
$ws.Range("AG10").Value = 251
$ws.Range("AH10").Value = 11
$ws.Range("AI10").Value = 21
$ws.Range("AN10").Value = 3.75
$ws.Range("AO10").Value = 10
$ws.Range("AP10").Value = 21
$ws.Range("AQ10").Value = 34
$ws.Range("AR10").Value = 51
$ws.Range("AS10").Value = 151
$ws.Range("AT10").Value = 2.75
$ws.Range("AU10").Value = 8
$ws.Range("AW10").Value = 6
$ws.Range("AY10").Value = 29
$ws.Range("AZ10").Value = 81
$ws.Range("BB10").Value = 201

# Row 11
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.53

# Row 12
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8

# Row 13
$ws.Range("G13").Value = 1.33
$ws.Range("H13").Value = 4.5
$ws.Range("I13").Value = 10
$ws.Range("J13").Value = 1.83
$ws.Range("K13").Value = 2.4
$ws.Range("Q13").Value = 1.85
$ws.Range("R13").Value = 2
$ws.Range("W13").Value = 6.5
$ws.Range("Z13").Value = 8
$ws.Range("AA13").Value = 12
$ws.Range("AC13").Value = 10
$ws.Range("AD13").Value = 9
$ws.Range("AE13").Value = 23
$ws.Range("AF13").Value = 81
$ws.Range("AH13").Value = 19
$ws.Range("AJ13").Value = 26
$ws.Range("AW13").Value = 9.5
